# Update "想去人数" (want-to-go count) figures across the workbook sheets.
# This corresponds to the gh-pages data refresh at commit 456a3b4.

$wb = $excel.ActiveWorkbook

# --- Sheet "展览" ---
$ws = $wb.Worksheets.Item("展览")
$ws.Range("F2").Value = 47
$ws.Range("F5").Value = 343
$ws.Range("F7").Value = 889
$ws.Range("F8").Value = 63
$ws.Range("F9").Value = 537
$ws.Range("F12").Value = 1165
$ws.Range("F14").Value = 252
$ws.Range("F17").Value = 6712
$ws.Range("F19").Value = 73
$ws.Range("F20").Value = 22
$ws.Range("F21").Value = 7621
$ws.Range("F23").Value = 37
$ws.Range("F24").Value = 3415
$ws.Range("F25").Value = 32
$ws.Range("F26").Value = 2142
$ws.Range("F27").Value = 913
$ws.Range("F29").Value = 167
$ws.Range("F32").Value = 1
$ws.Range("F34").Value = 201
$ws.Range("F35").Value = 1751
$ws.Range("F37").Value = 193
$ws.Range("F39").Value = 3
$ws.Range("F41").Value = 1238
$ws.Range("F42").Value = 1841
$ws.Range("F43").Value = 2148

# --- Sheet "本地生活" ---
$ws = $wb.Worksheets.Item("本地生活")
$ws.Range("F3").Value = 1236

# --- Sheet "全部类型" ---
$ws = $wb.Worksheets.Item("全部类型")
$ws.Range("F3").Value = 47
$ws.Range("F4").Value = 1236
$ws.Range("F7").Value = 343
$ws.Range("F9").Value = 889
$ws.Range("F10").Value = 63
$ws.Range("F11").Value = 537
$ws.Range("F14").Value = 1165
$ws.Range("F17").Value = 252
$ws.Range("F20").Value = 6712
$ws.Range("F22").Value = 73
$ws.Range("F23").Value = 22
$ws.Range("F24").Value = 7621
$ws.Range("F26").Value = 37
$ws.Range("F27").Value = 3415
$ws.Range("F28").Value = 32
$ws.Range("F29").Value = 2142
$ws.Range("F30").Value = 913
$ws.Range("F32").Value = 167
$ws.Range("F36").Value = 1
$ws.Range("F38").Value = 1751
$ws.Range("F40").Value = 193
$ws.Range("F42").Value = 3
$ws.Range("F44").Value = 1238
$ws.Range("F45").Value = 1841
$ws.Range("F47").Value = 2148
